$d = $word.ActiveDocument

# 1. First paragraph: "update the edited Madden NFL '08 Roster file ('base.ros')"
#    becomes "update the working base (already edited) Madden NFL '08 Roster file 'base.ros'"
$d.Content.Find.Execute(
    "update the edited Madden NFL ’08 Roster file (‘base.ros’) with the latest",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "update the working base (already edited) Madden NFL ’08 Roster file ‘base.ros’ with the latest",
    2) | Out-Null

# 2. Remove the trailing space inside "...EA ratings\edited" " before the period.
$d.Content.Find.Execute(
    "docs\EA ratings\edited” . Then make a copy",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "docs\EA ratings\edited”. Then make a copy",
    2) | Out-Null

# 3. "overwriting the previous copy there." -> "overwriting any previous copy there."
$d.Content.Find.Execute(
    ", overwriting the previous copy there.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    ", overwriting any previous copy there.",
    2) | Out-Null

# 4. "If" / " the format of the HTML on the NFL site has not changed" / ", the step 3
#    script should have created a file "" were three separate runs; combine them into one
#    (text itself is unchanged, so re-assert it to normalise the run boundaries).
$d.Content.Find.Execute(
    "If the format of the HTML on the NFL site has not changed, the step 3 script should have created a file “",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "If the format of the HTML on the NFL site has not changed, the step 3 script should have created a file “",
    2) | Out-Null

# 5. Move the "_GoBack" bookmark from the paragraph ending in "Hair Style, etc."
#    to the final (empty) paragraph of the document.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $lastPara.Range) | Out-Null
